$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1 — same text style (bold/border/centered) as the other
# header cells in row 1, so copy formatting from E1 (the previous last header).
$ws.Range("F1").Value = "Individual_Gain"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data column F2:F6 — same "Individual_Gain" value repeated for every run,
# matching Player_Life (column C) in this dataset.
$val = 87.40000000000012
$ws.Range("F2:F6").Value = $val
